$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price (D) cells that contain plain decimal numbers to remain
# stored as text, matching the original inline-string cell type, before
# assigning their new values (otherwise Excel would coerce them to numbers).
$textForceRefs = @("D4", "D5", "D6", "D7", "D8", "D9", "D11", "D12", "D13", "D14", "D15", "D16", "D17", "D19", "D20", "D21", "D22", "D23", "D24", "D26", "D27", "D28", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D50", "D51")
foreach ($ref in $textForceRefs) {
    $ws.Range($ref).NumberFormat = "@"
}

# Apply the updated coin data (column B/C swaps for re-ranked coins, plus
# refreshed Price and Volume(1h) figures for every row).
$ws.Range("D2").Value = "21.315.31"
$ws.Range("E2").Value = "  +3.89%  "
$ws.Range("D3").Value = "1.545.67"
$ws.Range("E3").Value = "  +4.76%  "
$ws.Range("D4").Value = "1.009"
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").Value = "0.9698"
$ws.Range("E5").Value = "  +0.76%  "
$ws.Range("D6").Value = "282.12"
$ws.Range("E6").Value = "  +2.05%  "
$ws.Range("D7").Value = "0.3632"
$ws.Range("E7").Value = "  -0.47%  "
$ws.Range("D8").Value = "0.3200"
$ws.Range("E8").Value = "  +4.40%  "
$ws.Range("D9").Value = "41.06"
$ws.Range("E9").Value = "  +3.16%  "
$ws.Range("E10").Value = "  +6.30%  "
$ws.Range("D11").Value = "0.06955"
$ws.Range("E11").Value = "  +5.31%  "
$ws.Range("D12").Value = "1.004"
$ws.Range("E12").Value = "  +0.27%  "
$ws.Range("B13").Value = "Solana"
$ws.Range("C13").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D13").Value = "19.04"
$ws.Range("E13").Value = "  +4.64%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "5.727"
$ws.Range("E14").Value = "  +4.94%  "
$ws.Range("D15").Value = "6.431"
$ws.Range("E15").Value = "  +4.29%  "
$ws.Range("D16").Value = "0.00001056"
$ws.Range("E16").Value = "  +2.81%  "
$ws.Range("D17").Value = "0.9696"
$ws.Range("E17").Value = "  +0.07%  "
$ws.Range("D18").Value = "1.547.52"
$ws.Range("E18").Value = "  +4.78%  "
$ws.Range("D19").Value = "0.06147"
$ws.Range("E19").Value = "  +4.30%  "
$ws.Range("D20").Value = "73.05"
$ws.Range("E20").Value = "  +6.15%  "
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").Value = "5.752"
$ws.Range("E21").Value = "  +5.27%  "
$ws.Range("B22").Value = "Avalanche"
$ws.Range("C22").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D22").Value = "15.36"
$ws.Range("E22").Value = "  +6.48%  "
$ws.Range("D23").Value = "11.45"
$ws.Range("E23").Value = "  +4.56%  "
$ws.Range("D24").Value = "2.324"
$ws.Range("E24").Value = "  +3.41%  "
$ws.Range("D25").Value = "21.375.81"
$ws.Range("E25").Value = "  +4.03%  "
$ws.Range("D26").Value = "147.75"
$ws.Range("E26").Value = "  +4.36%  "
$ws.Range("D27").Value = "2.297"
$ws.Range("E27").Value = "  +7.75%  "
$ws.Range("D28").Value = "17.88"
$ws.Range("E28").Value = "  +4.29%  "
$ws.Range("D29").Value = "1.719.90"
$ws.Range("E29").Value = "  +5.45%  "
$ws.Range("D30").Value = "119.09"
$ws.Range("E30").Value = "  +5.17%  "
$ws.Range("D31").Value = "4.041"
$ws.Range("E31").Value = "  +3.65%  "
$ws.Range("D32").Value = "0.8857"
$ws.Range("E32").Value = "  +9.37%  "
$ws.Range("D33").Value = "5.262"
$ws.Range("E33").Value = "  +6.24%  "
$ws.Range("D34").Value = "0.08067"
$ws.Range("E34").Value = "  +2.55%  "
$ws.Range("D35").Value = "1.562"
$ws.Range("E35").Value = "  +2.71%  "
$ws.Range("D36").Value = "1.225"
$ws.Range("E36").Value = "  -3.24%  "
$ws.Range("D37").Value = "5.011"
$ws.Range("E37").Value = "  +5.80%  "
$ws.Range("D38").Value = "0.05903"
$ws.Range("E38").Value = "  +3.07%  "
$ws.Range("D39").Value = "0.2004"
$ws.Range("E39").Value = "  +6.87%  "
$ws.Range("D40").Value = "0.02135"
$ws.Range("E40").Value = "  +4.79%  "
$ws.Range("D41").Value = "10.90"
$ws.Range("E41").Value = "  +4.82%  "
$ws.Range("D42").Value = "7.938"
$ws.Range("E42").Value = "  +3.51%  "
$ws.Range("D43").Value = "0.9691"
$ws.Range("E43").Value = "  +0.66%  "
$ws.Range("D44").Value = "0.5526"
$ws.Range("E44").Value = "  +4.76%  "
$ws.Range("D45").Value = "12.65"
$ws.Range("E45").Value = "  +5.26%  "
$ws.Range("D46").Value = "3.584"
$ws.Range("E46").Value = "  +2.28%  "
$ws.Range("D47").Value = "0.5511"
$ws.Range("E47").Value = "  +6.61%  "
$ws.Range("D48").Value = "122.35"
$ws.Range("E48").Value = "  +4.78%  "
$ws.Range("E49").Value = "  +6.77%  "
$ws.Range("D50").Value = "0.06621"
$ws.Range("E50").Value = "  +2.70%  "
$ws.Range("D51").Value = "70.07"
$ws.Range("E51").Value = "  +4.73%  "
